$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose npc_name (column C) is "심학규" (Shim Hak-gyu) -- these NPC rows are
# reclassified from object_type "NPC" to "PLAYER" (the player-controlled character).
$playerRows = @(
    21,57,114,180,181,219,221,237,240,243,244,247,248,267,284,293,350,413,414,415,417,418,419,420,
    421,422,423,424,425,426,427,428,430,431,432,433,434,435,436,437,439,442,443,444,445,446,448,451,
    452,453,454,456,457,458,459,460,461,462,463,464,466,468,471,472,477,478,479,480,482,486,487,492,
    493,494,495,497,498,499,500,503,504,505,506,515,516,517,518,519,520,521,522,523,524,525,526,527,
    528,529,530,531,532,533,534,535,536,537,538,539,540,541,542,543,544,545,546,547,548,549,550,551,
    552,553,554,555,556,557,558,559,560,561,562,563,564,606
)

foreach ($r in $playerRows) {
    $ws.Cells.Item($r, 2).Value = "PLAYER"
}

# Turn on AutoFilter over the full used range of the sheet.
$ws.Range("A1:F606").AutoFilter() | Out-Null

# Applying AutoFilter defines the workbook-scoped, hidden built-in
# _FilterDatabase name pointing at the filtered range.
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=NPC_01!`$A`$1:`$F`$606")
$fd.Visible = $false

# Move the active selection to A2.
$ws.Range("A2").Select() | Out-Null
